$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.924.72'
$ws.Range("E2").Value = '  -0.69%  '

# Row 3
$ws.Range("D3").Value = '2.351.86'
$ws.Range("E3").Value = '  -1.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").Value = "'0.675"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.93%  '

# Row 6
$ws.Range("D6").Value = "'241.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '

# Row 7
$ws.Range("D7").Value = "'73.02"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.83%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("E9").Value = '  +0.58%  '

# Row 10
$ws.Range("E10").Value = '  -2.64%  '

# Row 11
$ws.Range("D11").Value = "'58.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.67%  '

# Row 12
$ws.Range("D12").Value = "'33.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.35%  '

# Row 13
$ws.Range("E13").Value = '  -0.06%  '

# Row 14
$ws.Range("E14").Value = '  -2.91%  '

# Row 15
$ws.Range("D15").Value = '2.702.72'
$ws.Range("E15").Value = '  -0.95%  '

# Row 16
$ws.Range("D16").Value = "'16.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.74%  '

# Row 17
$ws.Range("E17").Value = '  -1.92%  '

# Row 18
$ws.Range("D18").Value = '2.350.53'
$ws.Range("E18").Value = '  -0.86%  '

# Row 19
$ws.Range("D19").Value = '43.817.46'
$ws.Range("E19").Value = '  -1.16%  '

# Row 20
$ws.Range("D20").Value = "'0.0000104"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.29%  '

# Row 21
$ws.Range("D21").Value = "'6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '

# Row 22
$ws.Range("D22").Value = "'78.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.32%  '

# Row 23
$ws.Range("D23").Value = "'255.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.99%  '

# Row 24
$ws.Range("D24").Value = "'1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.97%  '

# Row 25
$ws.Range("E25").Value = '  -0.07%  '

# Row 26
$ws.Range("E26").Value = '  +0.82%  '

# Row 27
$ws.Range("D27").Value = "'2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '

# Row 28
$ws.Range("E28").Value = '  -2.53%  '

# Row 29
$ws.Range("D29").Value = "'2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.52%  '

# Row 30
$ws.Range("D30").Value = "'22.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.17%  '

# Row 31
$ws.Range("D31").Value = "'177.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.15%  '

# Row 32
$ws.Range("E32").Value = '  -1.72%  '

# Row 33
$ws.Range("E33").Value = '  +1.01%  '

# Row 34
$ws.Range("E34").Value = '  -1.68%  '

# Row 35
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = "'5.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.95%  '

# Row 36
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = "'5.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.08%  '

# Row 37
$ws.Range("D37").Value = "'3.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.37%  '

# Row 38
$ws.Range("D38").Value = "'6.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.22%  '

# Row 39
$ws.Range("E39").Value = '  -4.51%  '

# Row 40
$ws.Range("D40").Value = "'0.0277"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '

# Row 41
$ws.Range("D41").Value = "'67.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +26.63%  '

# Row 42
$ws.Range("E42").Value = '  +14.16%  '

# Row 43
$ws.Range("D43").Value = "'0.110"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.20%  '

# Row 44
$ws.Range("D44").Value = "'9.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.94%  '

# Row 45
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = "'0.202"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.02%  '

# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = "'18.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.42%  '

# Row 47
$ws.Range("D47").Value = "'2.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.90%  '

# Row 48
$ws.Range("E48").Value = '  -1.61%  '

# Row 49
$ws.Range("E49").Value = '  +0.03%  '

# Row 50
$ws.Range("D50").Value = "'99.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.78%  '

# Row 51
$ws.Range("E51").Value = '  -4.97%  '

Write-Host "Applied cryptos update"